$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old rows 99-137 (return/sell/scene4_lily/ending scenes) before rewriting
$ws.Range("A99:K137").ClearContents()

$ws.Range("D99").Value = 'setFlag'
$ws.Range("E99").Value = 'chitsii.arena.player.kain_soul_choice,0'

$ws.Range("D100").Value = 'modInvoke'
$ws.Range("E100").Value = 'complete_quest(06_2_zek_steal_soulgem)'
$ws.Range("F100").Value = 'pc'

$ws.Range("D101").Value = 'modInvoke'
$ws.Range("E101").Value = 'complete_quest(06_2_zek_steal_soulgem_return)'
$ws.Range("F101").Value = 'pc'

$ws.Range("D102").Value = 'end'

$ws.Range("A103").Value = 'sell'

$ws.Range("F104").Value = 'sukutsu_shady_merchant'
$ws.Range("H104").Value = 'zek_sell1'
$ws.Range("I104").Value = 'ふふ、素晴らしい！ これです、これこそが私が求めていた『合理的かつ冷酷な決断』だ！'
$ws.Range("J104").Value = 'ふふ、素晴らしい！ これです、これこそが私が求めていた『合理的かつ冷酷な決断』だ！'

$ws.Range("F105").Value = 'sukutsu_shady_merchant'
$ws.Range("H105").Value = 'zek_sell2'
$ws.Range("I105").Value = '友情を燃料にして、さらなる高みへ昇る……。あなたは、本物の怪物の素質がある。'
$ws.Range("J105").Value = '友情を燃料にして、さらなる高みへ昇る……。あなたは、本物の怪物の素質がある。'

$ws.Range("F106").Value = 'pc'
$ws.Range("H106").Value = 'narr_sell1'
$ws.Range("I106").Value = '（彼は懐から何かを取り出す。）'
$ws.Range("J106").Value = '（彼は懐から何かを取り出す。）'

$ws.Range("F107").Value = 'sukutsu_shady_merchant'
$ws.Range("H107").Value = 'zek_sell3'
$ws.Range("I107").Value = 'さあ、約束の報酬です。**小さなコイン15枚**と**プラチナコイン5枚**を、台帳に記録する手はずを整えておきましょう。それと、この『暗い印』を。'
$ws.Range("J107").Value = 'さあ、約束の報酬です。**小さなコイン15枚**と**プラチナコイン5枚**を、台帳に記録する手はずを整えておきましょう。それと、この『暗い印』を。'

$ws.Range("D108").Value = 'shake'

$ws.Range("F109").Value = 'sukutsu_shady_merchant'
$ws.Range("H109").Value = 'zek_sell4'
$ws.Range("I109").Value = 'これで、あなたは『魂を喰らう者』となりました。……では、良い演技を。彼に気づかれないよう、お気をつけて。'
$ws.Range("J109").Value = 'これで、あなたは『魂を喰らう者』となりました。……では、良い演技を。彼に気づかれないよう、お気をつけて。'

$ws.Range("F110").Value = 'pc'
$ws.Range("H110").Value = 'narr_sell2'
$ws.Range("I110").Value = '（ゼクは影の中へと消えていく。）'
$ws.Range("J110").Value = '（ゼクは影の中へと消えていく。）'

$ws.Range("D111").Value = 'eval'
$ws.Range("E111").Value = 'for(int i=0; i<15; i++) { EClass.pc.Pick(ThingGen.Create("medal")); } for(int i=0; i<5; i++) { EClass.pc.Pick(ThingGen.Create("plat")); }'

$ws.Range("B112").Value = 'sell_balgas'

$ws.Range("A113").Value = 'sell_balgas'

$ws.Range("D114").Value = 'eval'
$ws.Range("E114").Value = 'Debug.Log("[SukutsuArena] Attempting to play BGM: BGM/Lobby_Normal");             var data = SoundManager.current.GetData("BGM/Lobby_Normal");             if (data != null) {                 Debug.Log("[SukutsuArena] Found BGM data, type: " + data.GetType().Name);                 if (data is BGMData bgm) {                     Debug.Log("[SukutsuArena] Playing as BGM");                     LayerDrama.haltPlaylist = true;                     LayerDrama.maxBGMVolume = true;                     SoundManager.current.PlayBGM(bgm);                 } else {                     Debug.Log("[SukutsuArena] Playing as Sound");                     SoundManager.current.Play(data);                 }             } else {                 Debug.LogWarning("[SukutsuArena] BGM not found: BGM/Lobby_Normal");             }'

$ws.Range("F115").Value = 'pc'
$ws.Range("H115").Value = 'narr_sell3'
$ws.Range("I115").Value = '（あなたはロビーに戻る。バルガスがあなたを待っている。）'
$ws.Range("J115").Value = '（あなたはロビーに戻る。バルガスがあなたを待っている。）'

$ws.Range("D116").Value = 'wait'
$ws.Range("E116").Value = '0.3'

$ws.Range("D117").Value = 'focusChara'
$ws.Range("E117").Value = 'sukutsu_arena_master'

$ws.Range("D118").Value = 'wait'
$ws.Range("E118").Value = '0.5'

$ws.Range("F119").Value = 'sukutsu_arena_master'
$ws.Range("H119").Value = 'balgas_sell1'
$ws.Range("I119").Value = '……おい。カインの魂の欠片は見つかったか？'
$ws.Range("J119").Value = '……おい。カインの魂の欠片は見つかったか？'

$ws.Range("F120").Value = 'pc'
$ws.Range("H120").Value = 'narr_sell4'
$ws.Range("I120").Value = '（あなたは首を横に振る。）'
$ws.Range("J120").Value = '（あなたは首を横に振る。）'

$ws.Range("F121").Value = 'sukutsu_arena_master'
$ws.Range("H121").Value = 'balgas_sell2'
$ws.Range("I121").Value = '……そうか。見つからなかったか。'
$ws.Range("J121").Value = '……そうか。見つからなかったか。'

$ws.Range("F122").Value = 'pc'
$ws.Range("H122").Value = 'narr_sell5'
$ws.Range("I122").Value = '（彼は深く息を吐き、酒瓶を手に取る。）'
$ws.Range("J122").Value = '（彼は深く息を吐き、酒瓶を手に取る。）'

$ws.Range("F123").Value = 'sukutsu_arena_master'
$ws.Range("H123").Value = 'balgas_sell3'
$ws.Range("I123").Value = '……まあ、仕方ねえ。お前は十分頑張った。……ありがよ。'
$ws.Range("J123").Value = '……まあ、仕方ねえ。お前は十分頑張った。……ありがよ。'

$ws.Range("D124").Value = 'setFlag'
$ws.Range("E124").Value = 'chitsii.arena.player.kain_soul_choice,1'

$ws.Range("D125").Value = 'modInvoke'
$ws.Range("E125").Value = 'complete_quest(06_2_zek_steal_soulgem)'
$ws.Range("F125").Value = 'pc'

$ws.Range("D126").Value = 'modInvoke'
$ws.Range("E126").Value = 'complete_quest(06_2_zek_steal_soulgem_sell)'
$ws.Range("F126").Value = 'pc'

$ws.Range("D127").Value = 'end'
